# Refresh market-price-derived columns (currentAveragePrice* / LevePrice* /
# LeveProfit*, columns H:N) on the affected leve rows across sheets, as
# produced by the scheduled market-data refresh runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 86
$ws.Cells.Item(86, 8).Value = 6237.273
$ws.Cells.Item(86, 9).Value = 6001.5
$ws.Cells.Item(86, 10).Value = 6372
$ws.Cells.Item(86, 11).Value = 6001.5
$ws.Cells.Item(86, 12).Value = 6372
$ws.Cells.Item(86, 13).Value = -4878.5
$ws.Cells.Item(86, 14).Value = -8618
# Row 89
$ws.Cells.Item(89, 8).Value = 6237.273
$ws.Cells.Item(89, 9).Value = 6001.5
$ws.Cells.Item(89, 10).Value = 6372
$ws.Cells.Item(89, 11).Value = 30007.5
$ws.Cells.Item(89, 12).Value = 31860
$ws.Cells.Item(89, 13).Value = -24391.5
$ws.Cells.Item(89, 14).Value = -43092

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value = 7203.3647
$ws.Cells.Item(32, 9).Value = 5011.8887
$ws.Cells.Item(32, 10).Value = 19754.545
$ws.Cells.Item(32, 11).Value = 5011.8887
$ws.Cells.Item(32, 12).Value = 19754.545
$ws.Cells.Item(32, 13).Value = -4724.8887
$ws.Cells.Item(32, 14).Value = -20328.545
# Row 63
$ws.Cells.Item(63, 8).Value = 3107
$ws.Cells.Item(63, 9).Value = 2933.3333
$ws.Cells.Item(63, 10).Value = 3211.2
$ws.Cells.Item(63, 11).Value = 2933.3333
$ws.Cells.Item(63, 12).Value = 3211.2
$ws.Cells.Item(63, 13).Value = -2247.3333
$ws.Cells.Item(63, 14).Value = -4583.2
# Row 66
$ws.Cells.Item(66, 8).Value = 3107
$ws.Cells.Item(66, 9).Value = 2933.3333
$ws.Cells.Item(66, 10).Value = 3211.2
$ws.Cells.Item(66, 11).Value = 14666.6665
$ws.Cells.Item(66, 12).Value = 16056
$ws.Cells.Item(66, 13).Value = -11234.6665
$ws.Cells.Item(66, 14).Value = -22920

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value = 1493.4642
$ws.Cells.Item(31, 9).Value = 1185.8636
$ws.Cells.Item(31, 10).Value = 2621.3333
$ws.Cells.Item(31, 11).Value = 1185.8636
$ws.Cells.Item(31, 12).Value = 2621.3333
$ws.Cells.Item(31, 13).Value = -890.8635999999999
$ws.Cells.Item(31, 14).Value = -3211.3333
# Row 34
$ws.Cells.Item(34, 8).Value = 1493.4642
$ws.Cells.Item(34, 9).Value = 1185.8636
$ws.Cells.Item(34, 10).Value = 2621.3333
$ws.Cells.Item(34, 11).Value = 1185.8636
$ws.Cells.Item(34, 12).Value = 2621.3333
$ws.Cells.Item(34, 13).Value = -983.8635999999999
$ws.Cells.Item(34, 14).Value = -3025.3333
# Row 44
$ws.Cells.Item(44, 8).Value = 3000
$ws.Cells.Item(44, 9).Value = 3000
$ws.Cells.Item(44, 11).Value = 3000
$ws.Cells.Item(44, 13).Value = -2558
# Row 62
$ws.Cells.Item(62, 8).Value = 2328
$ws.Cells.Item(62, 9).Value = 2344.4443
$ws.Cells.Item(62, 10).Value = 2285.7144
$ws.Cells.Item(62, 11).Value = 2344.4443
$ws.Cells.Item(62, 12).Value = 2285.7144
$ws.Cells.Item(62, 13).Value = -1720.4443
$ws.Cells.Item(62, 14).Value = -3533.7144
# Row 65
$ws.Cells.Item(65, 8).Value = 2328
$ws.Cells.Item(65, 9).Value = 2344.4443
$ws.Cells.Item(65, 10).Value = 2285.7144
$ws.Cells.Item(65, 11).Value = 11722.2215
$ws.Cells.Item(65, 12).Value = 11428.572
$ws.Cells.Item(65, 13).Value = -8602.2215
$ws.Cells.Item(65, 14).Value = -17668.572
# Row 141
$ws.Cells.Item(141, 8).Value = 32071.555
$ws.Cells.Item(141, 10).Value = 32071.555
$ws.Cells.Item(141, 12).Value = 32071.555
$ws.Cells.Item(141, 14).Value = -42431.555

$ws = $wb.Worksheets.Item("CUL")
# Row 33
$ws.Cells.Item(33, 8).Value = 141.25
$ws.Cells.Item(33, 9).Value = 140
$ws.Cells.Item(33, 10).Value = 141.66667
$ws.Cells.Item(33, 11).Value = 840
$ws.Cells.Item(33, 12).Value = 850.0000200000001
$ws.Cells.Item(33, 13).Value = -557
$ws.Cells.Item(33, 14).Value = -1416.00002
# Row 40
$ws.Cells.Item(40, 8).Value = 319.3684
$ws.Cells.Item(40, 9).Value = 83.454544
$ws.Cells.Item(40, 10).Value = 643.75
$ws.Cells.Item(40, 11).Value = 333.818176
$ws.Cells.Item(40, 12).Value = 2575
$ws.Cells.Item(40, 13).Value = -264.818176
$ws.Cells.Item(40, 14).Value = -2713
# Row 44
$ws.Cells.Item(44, 8).Value = 1022.2222
$ws.Cells.Item(44, 9).Value = 700
$ws.Cells.Item(44, 10).Value = 1280
$ws.Cells.Item(44, 11).Value = 2100
$ws.Cells.Item(44, 12).Value = 3840
$ws.Cells.Item(44, 13).Value = -1702
$ws.Cells.Item(44, 14).Value = -4636
# Row 68
$ws.Cells.Item(68, 8).Value = 983.3333
$ws.Cells.Item(68, 9).Value = 725
$ws.Cells.Item(68, 10).Value = 1500
$ws.Cells.Item(68, 11).Value = 2175
$ws.Cells.Item(68, 12).Value = 4500
$ws.Cells.Item(68, 13).Value = -1364
$ws.Cells.Item(68, 14).Value = -6122
# Row 69
$ws.Cells.Item(69, 8).Value = 5994
$ws.Cells.Item(69, 10).Value = 5994
$ws.Cells.Item(69, 12).Value = 17982
$ws.Cells.Item(69, 14).Value = -19604
# Row 71
$ws.Cells.Item(71, 8).Value = 983.3333
$ws.Cells.Item(71, 9).Value = 725
$ws.Cells.Item(71, 10).Value = 1500
$ws.Cells.Item(71, 11).Value = 6525
$ws.Cells.Item(71, 12).Value = 13500
$ws.Cells.Item(71, 13).Value = -2469
$ws.Cells.Item(71, 14).Value = -21612
# Row 72
$ws.Cells.Item(72, 8).Value = 5994
$ws.Cells.Item(72, 10).Value = 5994
$ws.Cells.Item(72, 12).Value = 53946
$ws.Cells.Item(72, 14).Value = -62058
# Row 80
$ws.Cells.Item(80, 8).Value = 3550.7273
$ws.Cells.Item(80, 9).Value = 0
$ws.Cells.Item(80, 10).Value = 3550.7273
$ws.Cells.Item(80, 11).Value = 0
$ws.Cells.Item(80, 12).Value = 10652.1819
$ws.Cells.Item(80, 13).ClearContents()
$ws.Cells.Item(80, 14).Value = -12524.1819
# Row 83
$ws.Cells.Item(83, 8).Value = 3550.7273
$ws.Cells.Item(83, 9).Value = 0
$ws.Cells.Item(83, 10).Value = 3550.7273
$ws.Cells.Item(83, 11).Value = 0
$ws.Cells.Item(83, 12).Value = 31956.5457
$ws.Cells.Item(83, 13).ClearContents()
$ws.Cells.Item(83, 14).Value = -41316.5457
# Row 86
$ws.Cells.Item(86, 8).Value = 453.86667
$ws.Cells.Item(86, 9).Value = 451
$ws.Cells.Item(86, 10).Value = 454.30768
$ws.Cells.Item(86, 11).Value = 1353
$ws.Cells.Item(86, 12).Value = 1362.92304
$ws.Cells.Item(86, 13).Value = -167
$ws.Cells.Item(86, 14).Value = -3734.92304
# Row 89
$ws.Cells.Item(89, 8).Value = 453.86667
$ws.Cells.Item(89, 9).Value = 451
$ws.Cells.Item(89, 10).Value = 454.30768
$ws.Cells.Item(89, 11).Value = 4059
$ws.Cells.Item(89, 12).Value = 4088.76912
$ws.Cells.Item(89, 13).Value = 1869
$ws.Cells.Item(89, 14).Value = -15944.76912
# Row 131
$ws.Cells.Item(131, 8).Value = 1071.2933
$ws.Cells.Item(131, 10).Value = 1125.3188
$ws.Cells.Item(131, 12).Value = 3375.9564
$ws.Cells.Item(131, 14).Value = -13455.9564

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Cells.Item(80, 8).Value = 2271.8635
$ws.Cells.Item(80, 9).Value = 2319.6875
$ws.Cells.Item(80, 10).Value = 2144.3333
$ws.Cells.Item(80, 11).Value = 2319.6875
$ws.Cells.Item(80, 12).Value = 2144.3333
$ws.Cells.Item(80, 13).Value = -1321.6875
$ws.Cells.Item(80, 14).Value = -4140.3333
# Row 83
$ws.Cells.Item(83, 8).Value = 2271.8635
$ws.Cells.Item(83, 9).Value = 2319.6875
$ws.Cells.Item(83, 10).Value = 2144.3333
$ws.Cells.Item(83, 11).Value = 11598.4375
$ws.Cells.Item(83, 12).Value = 10721.6665
$ws.Cells.Item(83, 13).Value = -6606.4375
$ws.Cells.Item(83, 14).Value = -20705.6665

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Cells.Item(7, 8).Value = 1664.5294
$ws.Cells.Item(7, 9).Value = 1706.1333
$ws.Cells.Item(7, 10).Value = 1352.5
$ws.Cells.Item(7, 11).Value = 1706.1333
$ws.Cells.Item(7, 12).Value = 1352.5
$ws.Cells.Item(7, 13).Value = -1594.1333
$ws.Cells.Item(7, 14).Value = -1576.5
# Row 68
$ws.Cells.Item(68, 8).Value = 9782.923000000001
$ws.Cells.Item(68, 9).Value = 34996.668
$ws.Cells.Item(68, 10).Value = 2218.8
$ws.Cells.Item(68, 11).Value = 34996.668
$ws.Cells.Item(68, 12).Value = 2218.8
$ws.Cells.Item(68, 13).Value = -34247.668
$ws.Cells.Item(68, 14).Value = -3716.8
# Row 71
$ws.Cells.Item(71, 8).Value = 9782.923000000001
$ws.Cells.Item(71, 9).Value = 34996.668
$ws.Cells.Item(71, 10).Value = 2218.8
$ws.Cells.Item(71, 11).Value = 174983.34
$ws.Cells.Item(71, 12).Value = 11094
$ws.Cells.Item(71, 13).Value = -171239.34
$ws.Cells.Item(71, 14).Value = -18582
# Row 126
$ws.Cells.Item(126, 8).Value = 1664.5294
$ws.Cells.Item(126, 9).Value = 1706.1333
$ws.Cells.Item(126, 10).Value = 1352.5
$ws.Cells.Item(126, 11).Value = 5118.3999
$ws.Cells.Item(126, 12).Value = 4057.5
$ws.Cells.Item(126, 13).Value = -2648.3999
$ws.Cells.Item(126, 14).Value = -8997.5

Write-Host "Applied $([int](188)) cell updates across 6 sheets."
